# feat: add 2022-Q3 data
#
# 1. Insert a new row at the top of the "总计" (summary) sheet's data table
#    holding the 2022-Q3 totals, shifting the existing quarters down by one
#    row and re-numbering the running index in column A.
# 2. Insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    before the former-first-quarter sheet "2022-Q2") with the per-fund
#    holdings detail for that quarter.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (leading apostrophe forces Excel to
# keep numeric-looking strings - fund codes, percentages, etc. - as text
# instead of silently coercing them to numbers).
function Set-TextCell($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q3 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# The inserted row picks up formatting copied from the row above (the
# header); clear that and re-apply the plain "index" style (column A only)
# used by every other data row so the new row matches its siblings.
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Application.CutCopyMode = $false

$summary.Cells.Item(2, 1).Value = 0
Set-TextCell $summary 2 2 "2022-Q3"
$summary.Cells.Item(2, 3).Value = 18
$summary.Cells.Item(2, 4).Value = 0.73

# The running index in column A (rows 3..8) needs to shift by +1 now that
# a new row sits above them.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6

# ---------------------------------------------------------------------
# Step 2: add the new "2022-Q3" detail sheet, right before "2022-Q2"
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Header row
Set-TextCell $q3Sheet 1 2 "基金代码"
Set-TextCell $q3Sheet 1 3 "基金名称"
Set-TextCell $q3Sheet 1 4 "基金规模"
Set-TextCell $q3Sheet 1 5 "股票总仓位"
Set-TextCell $q3Sheet 1 6 "仓位占比"
Set-TextCell $q3Sheet 1 7 "持有市值(亿元)"
Set-TextCell $q3Sheet 1 8 "仓位排名"

$q3data = @(
    @("519183", "万家双引擎灵活配置混合", "2.91", "94.03", "5.19", "0.1510", 10),
    @("009611", "兴全汇享一年持有期混合A", "10.95", "26.52", "1.19", "0.1303", 6),
    @("770001", "德邦优化灵活配置混合", "2.52", "92.61", "4.95", "0.1247", 2),
    @("001901", "前海开源沪港深隆鑫灵活配置混合A", "4.36", "38.21", "2.13", "0.0929", 7),
    @("005944", "工银聚福混合C", "4.28", "29.01", "1.36", "0.0582", 9),
    @("519097", "新华中小市值优选混合", "0.71", "67.35", "4.49", "0.0319", 3),
    @("002000", "工银新生利混合", "1.08", "28.83", "2.57", "0.0278", 5),
    @("003132", "德邦新回报灵活配置混合", "0.62", "72.03", "3.77", "0.0234", 4),
    @("009612", "兴全汇享一年持有期混合C", "1.88", "26.52", "1.19", "0.0224", 6),
    @("012977", "瑞达鑫红量化6个月持有混合A", "0.43", "94.69", "4.93", "0.0212", 5),
    @("005855", "中科沃土沃瑞灵活配置混合A", "0.71", "83.31", "2.75", "0.0195", 9),
    @("005856", "中科沃土沃瑞灵活配置混合C", "0.37", "83.31", "2.75", "0.0102", 9),
    @("004937", "中航混改精选混合C", "0.07", "81.93", "8.64", "0.0060", 7),
    @("012978", "瑞达鑫红量化6个月持有混合C", "0.11", "94.69", "4.93", "0.0054", 5),
    @("519099", "新华灵活主题混合", "0.15", "81.48", "2.12", "0.0032", 10),
    @("004936", "中航混改精选混合A", "0.01", "81.93", "8.64", "0.0009", 7),
    @("001902", "前海开源沪港深隆鑫灵活配置混合C", "0.04", "38.21", "2.13", "0.0009", 7),
    @("005943", "工银聚福混合A", "0.06", "29.01", "1.36", "0.0008", 9)
)

$r = 2
foreach ($row in $q3data) {
    $q3Sheet.Cells.Item($r, 1).Value = $r - 2
    Set-TextCell $q3Sheet $r 2 $row[0]
    $q3Sheet.Cells.Item($r, 3).Value = $row[1]
    Set-TextCell $q3Sheet $r 4 $row[2]
    Set-TextCell $q3Sheet $r 5 $row[3]
    Set-TextCell $q3Sheet $r 6 $row[4]
    Set-TextCell $q3Sheet $r 7 $row[5]
    $q3Sheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
